# Insert two new weekly price rows at the top of the "Cebollín" data block
# (worksheet row 74), pushing all existing data rows down by two. This is
# the weekly update referenced by the commit message ("Fruta / hortaliza,
# semanal"): two fresh quotes (Primera / Segunda quality, "$/paquete 36
# unidades") are added, and the two oldest rows that fall off the bottom
# of the originally-filled range are preserved as new rows at the end of
# the sheet (now 152 and 153), exactly as a normal row insertion would do.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows above row 74; this shifts rows 74:151 down to 76:153,
# automatically carrying the previous last two rows (150:151) down to the
# new last two rows (152:153). Row insertion carries down the formatting
# (e.g. the date number format on column D) from the row above, just like
# interactively inserting rows in Excel.
$insertRange = $ws.Range("A74:R75")
$insertRange.EntireRow.Insert()

# New row 74: Primera quality, $/paquete 36 unidades
$ws.Cells.Item(74, 1).Value = 7
$ws.Cells.Item(74, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(74, 3).Value = "Ñuble"
$ws.Cells.Item(74, 4).Value = 45090
$ws.Cells.Item(74, 5).Value = 16
$ws.Cells.Item(74, 6).Value = 100112037
$ws.Cells.Item(74, 7).Value = "Cebollín"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 150
$ws.Cells.Item(74, 11).Value = 6000
$ws.Cells.Item(74, 12).Value = 6000
$ws.Cells.Item(74, 13).Value = 6000
$ws.Cells.Item(74, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(74, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(74, 16).Value = 167
$ws.Cells.Item(74, 17).Value = 36
$ws.Cells.Item(74, 18).Value = "Hortaliza"

# New row 75: Segunda quality, $/paquete 36 unidades
$ws.Cells.Item(75, 1).Value = 7
$ws.Cells.Item(75, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(75, 3).Value = "Ñuble"
$ws.Cells.Item(75, 4).Value = 45090
$ws.Cells.Item(75, 5).Value = 16
$ws.Cells.Item(75, 6).Value = 100112037
$ws.Cells.Item(75, 7).Value = "Cebollín"
$ws.Cells.Item(75, 8).Value = "Sin especificar"
$ws.Cells.Item(75, 9).Value = "Segunda"
$ws.Cells.Item(75, 10).Value = 100
$ws.Cells.Item(75, 11).Value = 5000
$ws.Cells.Item(75, 12).Value = 5000
$ws.Cells.Item(75, 13).Value = 5000
$ws.Cells.Item(75, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(75, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(75, 16).Value = 139
$ws.Cells.Item(75, 17).Value = 36
$ws.Cells.Item(75, 18).Value = "Hortaliza"
